$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-12 Thursday" "2025-06-13 Friday"

Replace-Text "305×9=" "753×8="
Replace-Text "174×7=" "502×8="
Replace-Text "780×5=" "396×5="
Replace-Text "383×6=" "167×2="
Replace-Text "883×4=" "719×5="
Replace-Text "948×5=" "512×2="
Replace-Text "392×4=" "528×7="
Replace-Text "810×7=" "599×2="
Replace-Text "531×5=" "463×7="
Replace-Text "559×6=" "494×5="
Replace-Text "902×4=" "866×4="
Replace-Text "728×6=" "720×2="
Replace-Text "788×2=" "725×9="
Replace-Text "919×4=" "712×6="
Replace-Text "626×6=" "487×2="
Replace-Text "745×5=" "461×9="
Replace-Text "792×4=" "871×9="
Replace-Text "942×2=" "542×6="
Replace-Text "215×4=" "600×9="
Replace-Text "572×2=" "185×4="
Replace-Text "558×2=" "472×9="
Replace-Text "641×9=" "155×7="
Replace-Text "871×2=" "881×5="
Replace-Text "529×7=" "950×6="
Replace-Text "749×6=" "681×6="
